# Refactor CEF functions: rename e_cef_extrato to e_cef_extcef and
# e_cef_extratos to e_cef_extcefs; update references in the "Por empresa"
# lookup table.
#
# Concretely this means the "Pró soluto" row (row 4) on the "Por empresa"
# sheet no longer maps to the (now removed) "Informakon"/"rec" CEF
# extraction function, so its Origem/Arquivo/Tabela columns (C:E) are
# cleared. The sibling row (row 3, "Repasse PF da CEF") keeps referencing
# CEF/extcef as before.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Por empresa")

# Clear the now-obsolete CEF extraction reference (Origem/Arquivo/Tabela)
# for the "Pró soluto" row.
$ws.Range("C4:E4").ClearContents()

# Restore the last active cell selection as left by the author.
[void]$ws.Range("J25").Select()
